# Financial data update for KOF (Coca-Cola FEMSA) yearly financials worksheet.
# Updates numeric figures across Income Statement, Balance Sheet and Cash Flow
# Statement sections (columns D:J, the 7 yearly periods) to reflect refreshed
# source data, and flips several cells that no longer have data to "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KOF")

$ws.Range("D8").Value = 9478000
$ws.Range("E8").Value = 9191600
$ws.Range("F8").Value = 7880100
$ws.Range("G8").Value = 7618300
$ws.Range("H8").Value = 8068900
$ws.Range("I8").Value = 7641100
$ws.Range("J8").Value = 6373100

$ws.Range("D9").Value = 5159000
$ws.Range("E9").Value = 5071500
$ws.Range("F9").Value = 4154700
$ws.Range("G9").Value = 4081500
$ws.Range("H9").Value = 4296700
$ws.Range("I9").Value = 4091500
$ws.Range("J9").Value = 3449400

$ws.Range("D10").Value = 4319000
$ws.Range("E10").Value = 4120100
$ws.Range("F10").Value = 3725400
$ws.Range("G10").Value = 3536700
$ws.Range("H10").Value = 3772200
$ws.Range("I10").Value = 3549500
$ws.Range("J10").Value = 2923800

$ws.Range("D14").Value = 1302300
$ws.Range("I14").Value = 7700
$ws.Range("J14").Value = 6200

$ws.Range("D17").Value = 9782900
$ws.Range("E17").Value = 8129500
$ws.Range("F17").Value = 6734300
$ws.Range("G17").Value = 6519300
$ws.Range("H17").Value = 6968000
$ws.Range("I17").Value = 6533300
$ws.Range("J17").Value = 5441700

$ws.Range("D18").Value = -304900
$ws.Range("E18").Value = 1062100
$ws.Range("F18").Value = 1145800
$ws.Range("G18").Value = 1099000
$ws.Range("H18").Value = 1100900
$ws.Range("I18").Value = 1107700
$ws.Range("J18").Value = 931400

$ws.Range("D20").Value = 163700
$ws.Range("E20").Value = 71900
$ws.Range("F20").Value = -48400
$ws.Range("G20").Value = -45300
$ws.Range("H20").Value = -22300
$ws.Range("I20").Value = 36700
$ws.Range("J20").Value = 31000

$ws.Range("D21").Value = 461200
$ws.Range("E21").Value = 1581800
$ws.Range("F21").Value = 1466500
$ws.Range("G21").Value = 1412700
$ws.Range("H21").Value = 1447100
$ws.Range("I21").Value = 1438500
$ws.Range("J21").Value = "NA"

$ws.Range("D22").Value = 437900
$ws.Range("E22").Value = 386400
$ws.Range("F22").Value = 327700
$ws.Range("G22").Value = 286800
$ws.Range("H22").Value = 172800
$ws.Range("I22").Value = 101100
$ws.Range("J22").Value = 89400

$ws.Range("D23").Value = -579100
$ws.Range("E23").Value = 747600
$ws.Range("F23").Value = 769600
$ws.Range("G23").Value = 766900
$ws.Range("H23").Value = 905800
$ws.Range("I23").Value = 1043300
$ws.Range("J23").Value = 873000

$ws.Range("D24").Value = 216400
$ws.Range("E24").Value = 203200
$ws.Range("F24").Value = 235400
$ws.Range("G24").Value = 199700
$ws.Range("H24").Value = 296400
$ws.Range("I24").Value = 324500
$ws.Range("J24").Value = 293100

$ws.Range("D26").Value = -795500
$ws.Range("E26").Value = 544500
$ws.Range("F26").Value = 534200
$ws.Range("G26").Value = 567200
$ws.Range("H26").Value = 609400
$ws.Range("I26").Value = 718800
$ws.Range("J26").Value = 579900

$ws.Range("D27").Value = -851700
$ws.Range("E27").Value = 520800
$ws.Range("F27").Value = 529400
$ws.Range("G27").Value = 545200
$ws.Range("H27").Value = 597000
$ws.Range("I27").Value = 698900
$ws.Range("J27").Value = 576600

$ws.Range("D29").Value = 192700
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"

$ws.Range("D32").Value = -163700
$ws.Range("E32").Value = -71900
$ws.Range("F32").Value = 48400
$ws.Range("G32").Value = 45300
$ws.Range("H32").Value = 22300
$ws.Range("I32").Value = -36700
$ws.Range("J32").Value = -31000

$ws.Range("D33").Value = -659000
$ws.Range("E33").Value = 520800
$ws.Range("F33").Value = 529400
$ws.Range("G33").Value = 545200
$ws.Range("H33").Value = 597000
$ws.Range("I33").Value = 698900
$ws.Range("J33").Value = 576600

$ws.Range("D35").Value = -659000
$ws.Range("E35").Value = 520800
$ws.Range("F35").Value = 529400
$ws.Range("G35").Value = 545200
$ws.Range("H35").Value = 597000
$ws.Range("I35").Value = 698900
$ws.Range("J35").Value = 576600

$ws.Range("D41").Value = 970600
$ws.Range("E41").Value = 541800
$ws.Range("F41").Value = 827000
$ws.Range("G41").Value = 1340400
$ws.Range("H41").Value = 1583300
$ws.Range("I41").Value = 2116600
$ws.Range("J41").Value = 612500

$ws.Range("J42").Value = 17100

$ws.Range("D43").Value = 1176500
$ws.Range("E43").Value = 1002200
$ws.Range("F43").Value = 717200
$ws.Range("G43").Value = 1450100
$ws.Range("H43").Value = 1243100
$ws.Range("I43").Value = 929500
$ws.Range("J43").Value = 561000

$ws.Range("D44").Value = 587700
$ws.Range("E44").Value = 555700
$ws.Range("F44").Value = 417200
$ws.Range("G44").Value = 808800
$ws.Range("H44").Value = 944400
$ws.Range("I44").Value = 513000
$ws.Range("J44").Value = 385000

$ws.Range("D45").Value = 143700
$ws.Range("E45").Value = 251100
$ws.Range("F45").Value = 222900
$ws.Range("G45").Value = 540300
$ws.Range("H45").Value = 701000
$ws.Range("I45").Value = 217700
$ws.Range("J45").Value = 112700

$ws.Range("D46").Value = 2878600
$ws.Range("E46").Value = 2350800
$ws.Range("F46").Value = 2184200
$ws.Range("G46").Value = 1972000
$ws.Range("H46").Value = 2235900
$ws.Range("I46").Value = 2373800
$ws.Range("J46").Value = 1692500

$ws.Range("D47").Value = 665200
$ws.Range("E47").Value = 1162400
$ws.Range("F47").Value = 933400
$ws.Range("G47").Value = 1807800
$ws.Range("H47").Value = 1789700
$ws.Range("I47").Value = 321300
$ws.Range("J47").Value = 230700

$ws.Range("D48").Value = 3921800
$ws.Range("E48").Value = 3376700
$ws.Range("F48").Value = 2613500
$ws.Range("G48").Value = 5226500
$ws.Range("H48").Value = 5356600
$ws.Range("I48").Value = 4208200
$ws.Range("J48").Value = 1974700

$ws.Range("D49").Value = 6425800
$ws.Range("E49").Value = 6411400
$ws.Range("F49").Value = 4693800
$ws.Range("G49").Value = 15054200
$ws.Range("H49").Value = 10474900
$ws.Range("I49").Value = 6931800
$ws.Range("J49").Value = 6870400

$ws.Range("D52").Value = 883800
$ws.Range("E52").Value = 1141800
$ws.Range("F52").Value = 449100
$ws.Range("G52").Value = 955800
$ws.Range("H52").Value = 150100
$ws.Range("I52").Value = 424700
$ws.Range("J52").Value = 263400

$ws.Range("D54").Value = 14775200
$ws.Range("E54").Value = 14443100
$ws.Range("F54").Value = 10874100
$ws.Range("G54").Value = 10983600
$ws.Range("H54").Value = 11205900
$ws.Range("I54").Value = 8590800
$ws.Range("J54").Value = 7330700

$ws.Range("D57").Value = 1621600
$ws.Range("E57").Value = 1440100
$ws.Range("F57").Value = 1045500
$ws.Range("G57").Value = 1076000
$ws.Range("H57").Value = 1094900
$ws.Range("I57").Value = 1207500
$ws.Range("J57").Value = 856100

$ws.Range("D58").Value = 629500
$ws.Range("E58").Value = 157800
$ws.Range("F58").Value = 179500
$ws.Range("G58").Value = 62400
$ws.Range("H58").Value = 370900
$ws.Range("I58").Value = 482700
$ws.Range("J58").Value = 286500

$ws.Range("D59").Value = 624300
$ws.Range("E59").Value = 464000
$ws.Range("F59").Value = 351500
$ws.Range("G59").Value = 1660900
$ws.Range("H59").Value = 1046500
$ws.Range("I59").Value = 582100
$ws.Range("J59").Value = 243900

$ws.Range("D60").Value = 2875300
$ws.Range("E60").Value = 2062000
$ws.Range("F60").Value = 1576400
$ws.Range("G60").Value = 1469000
$ws.Range("H60").Value = 1675600
$ws.Range("I60").Value = 1528300
$ws.Range("J60").Value = 1333100

$ws.Range("D61").Value = 3681900
$ws.Range("E61").Value = 4440500
$ws.Range("F61").Value = 3282900
$ws.Range("G61").Value = 3361600
$ws.Range("H61").Value = 2941600
$ws.Range("I61").Value = 1295200
$ws.Range("J61").Value = 870000

$ws.Range("D62").Value = 940500
$ws.Range("E62").Value = 1256700
$ws.Range("F62").Value = 391000
$ws.Range("G62").Value = 588600
$ws.Range("H62").Value = 1051800
$ws.Range("I62").Value = 654400
$ws.Range("J62").Value = 533000

$ws.Range("D66").Value = 8435900
$ws.Range("E66").Value = 8126200
$ws.Range("F66").Value = 5456500
$ws.Range("G66").Value = 5515900
$ws.Range("H66").Value = 5355800
$ws.Range("I66").Value = 3333600
$ws.Range("J66").Value = 2674400

$ws.Range("D72").Value = 3195600
$ws.Range("E72").Value = 4219300
$ws.Range("F72").Value = 4057600
$ws.Range("G72").Value = 7719100
$ws.Range("H72").Value = 3625300
$ws.Range("I72").Value = 6672000
$ws.Range("J72").Value = 3486300

$ws.Range("D76").Value = 6339300
$ws.Range("E76").Value = 6316900
$ws.Range("F76").Value = 5417600
$ws.Range("G76").Value = 5467700
$ws.Range("H76").Value = 5850100
$ws.Range("I76").Value = 5257300
$ws.Range("J76").Value = 4656200

$ws.Range("D81").Value = -659000
$ws.Range("E81").Value = 520800
$ws.Range("F81").Value = 529400
$ws.Range("G81").Value = 545200
$ws.Range("H81").Value = 597000
$ws.Range("I81").Value = 698900
$ws.Range("J81").Value = 576600

$ws.Range("D83").Value = 602900
$ws.Range("E83").Value = 448200
$ws.Range("F83").Value = 369500
$ws.Range("G83").Value = 359400
$ws.Range("H83").Value = 368900
$ws.Range("I83").Value = 294400
$ws.Range("J83").Value = "NA"

$ws.Range("D89").Value = 1719000
$ws.Range("E89").Value = 1678100
$ws.Range("F89").Value = 1200000
$ws.Range("G89").Value = 1262300
$ws.Range("H89").Value = 1142900
$ws.Range("I89").Value = 1223200
$ws.Range("J89").Value = 718500

$ws.Range("D91").Value = -572500
$ws.Range("E91").Value = -533100
$ws.Range("F91").Value = -545400
$ws.Range("G91").Value = -561800
$ws.Range("H91").Value = -549000
$ws.Range("I91").Value = -503800
$ws.Range("J91").Value = -354500

$ws.Range("D94").Value = -563200
$ws.Range("E94").Value = -1392000
$ws.Range("F94").Value = -566100
$ws.Range("G94").Value = -576000
$ws.Range("H94").Value = -2559200
$ws.Range("I94").Value = -568400
$ws.Range("J94").Value = "NA"

$ws.Range("D96").Value = -361600
$ws.Range("E96").Value = -362700
$ws.Range("F96").Value = -331800
$ws.Range("G96").Value = -311900
$ws.Range("H96").Value = -310400
$ws.Range("I96").Value = -296500
$ws.Range("J96").Value = -225800

$ws.Range("D100").Value = -557300
$ws.Range("E100").Value = -503400
$ws.Range("F100").Value = -443100
$ws.Range("G100").Value = -587000
$ws.Range("H100").Value = 1215700
$ws.Range("I100").Value = 3100
$ws.Range("J100").Value = "NA"

$ws.Range("D101").Value = -169600
$ws.Range("E101").Value = -67800
$ws.Range("F101").Value = -34100
$ws.Range("G101").Value = -220700
$ws.Range("H101").Value = -208800
$ws.Range("I101").Value = -69400
$ws.Range("J101").Value = "NA"

$ws.Range("D102").Value = 428800
$ws.Range("E102").Value = -285100
$ws.Range("F102").Value = 156800
$ws.Range("G102").Value = -121400
$ws.Range("H102").Value = -409400
$ws.Range("I102").Value = 588500
$ws.Range("J102").Value = -15500
